# Update sample config:
#  - "meta" sheet: body-fields value changes from "[]" to "[data.field, data.something]"
#  - selection/active-cell bookkeeping: "pages" -> A3, "meta" -> B5 (meta stays the active tab)

$wb = $excel.ActiveWorkbook

# --- content change -------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("meta")
$wsMeta.Range("B4").Value = "[data.field, data.something]"

# --- view/selection changes ------------------------------------------------
$wsPages = $wb.Worksheets.Item("pages")
$wsPages.Activate()
$wsPages.Range("A3").Select()

# Re-activate "meta" last so it remains the selected tab, and update its
# active cell to B5.
$wsMeta.Activate()
$wsMeta.Range("B5").Select()
